# Add a new data row (row 85) at the end of the sheet, mirroring the
# existing JRG_* daily COVID tracking columns A:I.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(85, 1).Value = 207  # A85 - JRG_TOTAL
$ws.Cells.Item(85, 2).Value = 132  # B85 - JRG_RECUPERADOS
$ws.Cells.Item(85, 3).Value = 70   # C85 - JRG_ISOLAMENTO
$ws.Cells.Item(85, 4).Value = 2    # D85 - JRG_INTERNADOS
$ws.Cells.Item(85, 5).Value = 3    # E85 - JRG_MORTES
$ws.Cells.Item(85, 6).Value = 84   # F85 - JRG_N
$ws.Cells.Item(85, 7).Value = 72   # G85 - JRG_CASOS_ATIVOS
$ws.Cells.Item(85, 8).Value = 8    # H85 - JRG_CASOS_DIA
$ws.Cells.Item(85, 9).Value = 0    # I85 - JRG_MORTES_DIA
